$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Update the "datetimeFigureOut" field text (31/08/2020 -> 09/09/2020) on
#    the slide master and every slide layout. The field placeholder is
#    identified robustly via PlaceholderFormat.Type (16 = date placeholder).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq "31/08/2020") {
                    $tr.Text = "09/09/2020"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2. Slide 1, "TextBox 131": split the trailing run so the text reads
#    " - Aug 2020, CC " + "BY-SA 2.0" (was " - Aug 2020, CC BY-SA"), and move
#    / resize the textbox to its new position.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$caption = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 131") {
        $caption = $shp
    }
}

$tr = $caption.TextFrame.TextRange
$secondRun = $tr.Runs(2, 1)
$secondRun.Text = " - Aug 2020, CC "
$tr.Text = "https://cern.ch/Luca.Canali/docs/SparkExecutorMemory.png - Aug 2020, CC BY-SA 2.0"

$caption.Left = 465.91741943359375
$caption.Top = 518.2410278320312
$caption.Width = 478.9216613769531
$caption.Height = 23.02267837524414

# ---------------------------------------------------------------------------
# 3. Slide 1, "Rectangle 47": shrink its height (only the extent cy changes).
# ---------------------------------------------------------------------------
$rect = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 47") {
        $rect = $shp
    }
}

$rect.Height = 60.55937194824219
